$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jarno")

# --- Row 15: date + "0.25" (stored as text, matching the other 0.5/0.75
#     text entries elsewhere in the column) + description text ---

# Copy the date-formatted style (A13, numFmtId 14) onto A15:A16 first,
# then write the serial date values.
$ws.Cells.Item(13, 1).Copy() | Out-Null
$ws.Range("A15:A16").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15, 1).Value = 44986
$ws.Cells.Item(16, 1).Value = 44986

# B15 needs to hold the literal text "0.25" (not the number 0.25) while
# keeping its original style (s=12, General number format, no quote
# prefix) -- exactly like the existing "0.5"/"0.75" text cells elsewhere
# in this sheet. Temporarily mark the cell as Text so the value isn't
# auto-coerced to a number, then restore the original formatting (from a
# sibling cell that already carries the untouched style) via a
# formats-only paste so the style id is reused rather than a new one
# minted.
$ws.Cells.Item(15, 2).NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = "0.25"
$ws.Cells.Item(13, 2).Copy() | Out-Null
$ws.Cells.Item(15, 2).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(15, 3).Value = "Sprint-planning, Retroplanning"

# --- Row 16: date (already set above) + numeric 2.5 hours formatted as
#     "0.0" + description text ---
$ws.Cells.Item(16, 2).NumberFormat = "0.0"
$ws.Cells.Item(16, 2).Value = 2.5
$ws.Cells.Item(16, 3).Value = "Serviceview ja MainWindow"

# --- Rows 17-24: column B adopts the same "0.0" number style as B16,
#     cells stay empty ---
$ws.Range("B17:B24").NumberFormat = "0.0"

# --- Sheet activation / selection bookkeeping ---
$ws.Activate()
$ws.Range("C16").Select()
